# Weekly update: a new price observation is inserted as row 405 (pushing all
# subsequent rows down by one), extending the sheet from A1:R473 to A1:R474.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 405; Excel automatically shifts rows 405-473 down to
# 406-474 and grows the sheet dimension to A1:R474.
$ws.Rows(405).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(405, 1).Value = 4
$ws.Cells.Item(405, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(405, 3).Value = "Los Lagos"
$ws.Cells.Item(405, 4).Value = 45209
$ws.Cells.Item(405, 5).Value = 10
$ws.Cells.Item(405, 6).Value = 100112032
$ws.Cells.Item(405, 7).Value = "Zapallo italiano"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 220
$ws.Cells.Item(405, 11).Value = 26000
$ws.Cells.Item(405, 12).Value = 26000
$ws.Cells.Item(405, 13).Value = 26000
$ws.Cells.Item(405, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(405, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(405, 16).Value = 520
$ws.Cells.Item(405, 17).Value = 50
$ws.Cells.Item(405, 18).Value = "Hortaliza"
